$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[58.65068788193136, 67.64782710921654]"
$ws.Range("T2").Value = "[46.16938927672991, 52.36486293403507]"
$ws.Range("L3").Value = "[58.454422723617526, 69.69420129962971]"
$ws.Range("T3").Value = "[47.29396165767887, 52.962440801436365]"
